# EasyShip.xlsx - "Add files via upload" change:
# Populate row 2 of the EasyShip sheet with one shipment's order data.
# Cells are written in the same order the source workbook's shared-string
# table was built in, so new strings are appended in matching order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = "USD"
$ws.Range("E2").Value = "T924JII7D79U7RE"
$ws.Range("X2").Value = "LIC Football Merch"
$ws.Range("G2").Value = "Timmy Donley"
$ws.Range("I2").Value = "alexsilvestrini@yahoo.com"
$ws.Range("H2").Value = "+(516) 729-0379"
$ws.Range("K2").Value = "5 Bond Place"

# Postal code "11706" looks like a number, and a plain .Value assignment
# would store it as a numeric cell. The source file stores it as text
# (shared string) with no special cell style, so round it through a
# TEXT() formula and paste back as a value to force a genuine text cell
# without leaving a quote-prefix / number-format style behind.
$ws.Range("M2").Formula = "=TEXT(11706,""0"")"
$ws.Range("M2").Copy()
$ws.Range("M2").PasteSpecial(-4163)

$ws.Range("N2").Value = "Bayshore"

# Item Customs Value* is numeric.
$ws.Range("AA2").Value = 25.0
